$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 26, shifting existing rows 26:65 down to 27:66.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly price observation.
$ws.Range("A26").Value = 3
$ws.Range("B26").Value = "Femacal de La Calera"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44775
$ws.Range("E26").Value = 5
$ws.Range("F26").Value = 100112035
$ws.Range("G26").Value = "Bruselas (repollito)"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 93
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = 14516
$ws.Range("N26").Value = "$/malla 15 kilos"
$ws.Range("O26").Value = "Provincia de Quillota"
$ws.Range("P26").Value = 968
$ws.Range("Q26").Value = 15
$ws.Range("R26").Value = "Hortaliza"
